$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D to make room for the newest quarter (Sep 2018),
# shifting the existing quarterly columns one position to the right (D->E ... K->L).
$ws.Columns("D:D").Insert()

# Copy number formatting from column E (which now holds the formerly-column-D data)
# into the newly inserted column D so the new quarter's cells are formatted consistently
# (date format on the "Period Ending" rows, number format elsewhere).
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$rowData = @{
    7 = @(43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    8 = @(618800, 629600, 602000, 595300, 535400, 641500, 615200, 674400, 680700)
    9 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    10 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    15 = @(-34700, -31000, -29400, -29700, -30300, -30600, -29700, -15200, -16400)
    17 = @(393200, 411900, 389500, 472800, 387600, 457900, 449300, 545300, 511200)
    18 = @(225600, 217700, 212500, 122500, 147800, 183600, 165900, 129200, 169500)
    20 = @(-139800, -89700, -157600, -215400, -186500, -83200, -153100, -259100, -137000)
    21 = @(120500, 159000, 84300, -63200, -8500, 131000, 42500, -99400, 65100)
    22 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    23 = @(85800, 128000, 54900, -92900, -38800, 100400, 12800, -129900, 32400)
    24 = @(19900, 42300, -8400, -47500, -32600, 2800, -19700, -50600, 11200)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(65900, 85700, 63300, -45400, -6200, 97600, 32400, -79300, 21300)
    27 = @(63100, 85200, 62900, -40600, -5800, 95000, 35900, -78000, 29500)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    29 = @("NA", "NA", "NA", 0, "NA", "NA", "NA", -300, -400)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(139800, 89700, 157600, 215400, 186500, 83200, 153100, 259100, 137000)
    33 = @(63100, 85200, 62900, -40600, -5800, 95000, 35900, -78300, 29100)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(63100, 85200, 62900, -40600, -5800, 95000, 35900, -78300, 29100)
    38 = @(43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    41 = @(2013800, 1777900, 1610200, 786600, 2127600, 2151200, 1556100, 1175400, 1738200)
    42 = @(2990100, 2885200, 2928900, 3452100, 3513500, 3689500, 4470500, 4395100, 5280100)
    43 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    44 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    45 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    46 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    47 = @(15800, 15800, 15600, 15300, 32700, 29300, 29700, 30200, 25700)
    48 = @(134300, 183100, 184100, 192000, 208200, 186200, 181500, 182800, 188100)
    49 = @(2402300, 2407500, 2384800, 2359700, 2401400, 2399600, 2444100, 2503000, 2406400)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(239800, 243300, 249500, 236800, 521800, 473900, 465400, 433400, 355600)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(42740200, 43116300, 42238200, 41249300, 42832000, 42573500, 42597400, 43950300, 45520400)
    57 = @(1205500, 1469300, 1396500, 672600, 1199800, 1056700, 594000, 388900, 578200)
    58 = @(1825500, 1865000, 1190000, 1612000, 2367300, 2332600, 2246200, 1324600, "NA")
    59 = @(600, 1400, 21400, 2000, 1000, 1000, 1400, 400, 0)
    60 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    61 = @(10935400, 11932500, 11586700, 11006700, 10972600, 10951600, 10662200, 10815100, 11479200)
    62 = @(313100, 306300, 262400, 295700, 597700, 545300, 564400, 567500, 619100)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(37889700, 38308600, 37490900, 36560100, 38087200, 37817200, 37894700, 39158300, 40703600)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(2097000, 2052900, 1969500, 1852100, 1978500, 1983000, 1929500, 1949700, 1998300)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(4850400, 4807700, 4747400, 4689100, 4744700, 4756200, 4702700, 4792000, 4816900)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    80 = @(43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    81 = @(63100, 85200, 62900, -40600, -5800, 95000, 35900, -78300, 29100)
    83 = @(34700, 31000, 29400, 29700, 30300, 30600, 29700, 30600, 32600)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(394900, -68700, 566900, 163100, 50000, -647100, -1610600, -509200, -13500)
    91 = @(-29100, -32800, -23000, -25100, -48600, -34300, -20200, -44100, 29300)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-29000, -11700, -23000, -25100, -48600, -34300, -20200, -47800, -41600)
    96 = @(0, -18900, -14900, 0, 0, 0, -900, 0, 0)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(-300200, 169000, -375700, -371500, -399100, 197600, 1086900, 349500, 100800)
    101 = @(-800, 92300, -21700, 77800, 25700, 69000, -45000, 29000, -7000)
    102 = @(64900, 180900, 146600, -155700, -371900, -414700, -588900, -178500, 40500)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $col = 4
    foreach ($v in $vals) {
        if ($null -ne $v) {
            $ws.Cells.Item([int]$r, $col).Value = $v
        }
        $col = $col + 1
    }
}

Write-Host "Update complete"
